$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 1.01
$ws.Range("T2").Value = 1.04
$ws.Range("AB3").Value = 8.4
$ws.Range("AH3").Value = 28
$ws.Range("AK3").Value = 27
$ws.Range("AL3").Value = 60
$ws.Range("F3").Value = 1.77
$ws.Range("G3").Value = 1.82
$ws.Range("H3").Value = 5.1
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 3.4
$ws.Range("L3").Value = 1.47
$ws.Range("N3").Value = 3
$ws.Range("P3").Value = 1.68
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.25
$ws.Range("T3").Value = 2.02
$ws.Range("U3").Value = 1.8
$ws.Range("W3").Value = 1.96
$ws.Range("Y3").Value = 19
$ws.Range("AB4").Value = 12.5
$ws.Range("AC4").Value = 13
$ws.Range("AF4").Value = 10.5
$ws.Range("AG4").Value = 12
$ws.Range("AH4").Value = 27
$ws.Range("AK4").Value = 14.5
$ws.Range("AL4").Value = 32
$ws.Range("AN4").Value = 5.1
$ws.Range("F4").Value = 1.39
$ws.Range("G4").Value = 1.43
$ws.Range("H4").Value = 8.6
$ws.Range("I4").Value = 10.5
$ws.Range("J4").Value = 5.1
$ws.Range("K4").Value = 5.6
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 5.5
$ws.Range("P4").Value = 2.54
$ws.Range("Q4").Value = 1.52
$ws.Range("R4").Value = 1.63
$ws.Range("S4").Value = 2.28
$ws.Range("T4").Value = 1.79
$ws.Range("U4").Value = 2.02
$ws.Range("V4").Value = 1.11
$ws.Range("W4").Value = 3.25
$ws.Range("X4").Value = 32
$ws.Range("Y4").Value = 38
$ws.Range("Z4").Value = 110
$ws.Range("F5").Value = 3.6
$ws.Range("G5").Value = 4.2
$ws.Range("H5").Value = 1.98
$ws.Range("I5").Value = 2.18
$ws.Range("J5").Value = 3.15
$ws.Range("K5").Value = 4.2
$ws.Range("L5").Value = 1.32
$ws.Range("M5").Value = 1.06
$ws.Range("P5").Value = 1.97
$ws.Range("Q5").Value = 1.83
$ws.Range("S5").Value = 2.82
$ws.Range("T5").Value = 1.71
$ws.Range("W5").Value = 1.32
$ws.Range("AB6").Value = 38
$ws.Range("AC6").Value = 15
$ws.Range("AD6").Value = 11
$ws.Range("AG6").Value = 26
$ws.Range("AH6").Value = 20
$ws.Range("AI6").Value = 1000
$ws.Range("AK6").Value = 1000
$ws.Range("AM6").Value = 1000
$ws.Range("AO6").Value = 6.4
$ws.Range("F6").Value = 5.8
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 1.5
$ws.Range("I6").Value = 1.59
$ws.Range("J6").Value = 4.7
$ws.Range("K6").Value = 5.6
$ws.Range("L6").Value = 1.21
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 6
$ws.Range("O6").Value = 1.16
$ws.Range("P6").Value = 2.72
$ws.Range("Q6").Value = 1.48
$ws.Range("R6").Value = 1.7
$ws.Range("S6").Value = 2.16
$ws.Range("T6").Value = 1.63
$ws.Range("V6").Value = 2.68
$ws.Range("X6").Value = 1000
$ws.Range("Y6").Value = 16
$ws.Range("Z6").Value = 15
$ws.Range("AC7").Value = 9.800000000000001
$ws.Range("F7").Value = 3.7
$ws.Range("G7").Value = 4.7
$ws.Range("H7").Value = 1.94
$ws.Range("I7").Value = 2.18
$ws.Range("J7").Value = 3.1
$ws.Range("K7").Value = 4.2
$ws.Range("L7").Value = 1.45
$ws.Range("P7").Value = 1.74
$ws.Range("R7").Value = 1.28
$ws.Range("S7").Value = 3.4
$ws.Range("T7").Value = 1.87
$ws.Range("U7").Value = 1.89
$ws.Range("V7").Value = 1.84
$ws.Range("W7").Value = 1.29
$ws.Range("F8").Value = 2.34
$ws.Range("G8").Value = 2.86
$ws.Range("H8").Value = 2.62
$ws.Range("I8").Value = 3.1
$ws.Range("J8").Value = 3.6
$ws.Range("L8").Value = 1.31
$ws.Range("N8").Value = 4.1
$ws.Range("P8").Value = 2.24
$ws.Range("S8").Value = 2.4
$ws.Range("T8").Value = 1.62
$ws.Range("U8").Value = 2.38
$ws.Range("V8").Value = 1.48
$ws.Range("W8").Value = 1.57
$ws.Range("F9").Value = 1.51
$ws.Range("G9").Value = 1.61
$ws.Range("H9").Value = 5.8
$ws.Range("I9").Value = 8.6
$ws.Range("J9").Value = 4.1
$ws.Range("K9").Value = 5
$ws.Range("M9").Value = 1.06
$ws.Range("P9").Value = 1.96
$ws.Range("Q9").Value = 1.84
$ws.Range("R9").Value = 1.37
$ws.Range("U9").Value = 1.86
$ws.Range("W9").Value = 2.62
$ws.Range("AB10").Value = 14.5
$ws.Range("AC10").Value = 11.5
$ws.Range("AD10").Value = 19.5
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 14.5
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 18
$ws.Range("AJ10").Value = 22
$ws.Range("AK10").Value = 17
$ws.Range("AL10").Value = 26
$ws.Range("AM10").Value = 75
$ws.Range("AN10").Value = 8
$ws.Range("F10").Value = 1.69
$ws.Range("G10").Value = 1.71
$ws.Range("H10").Value = 4.4
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 4.6
$ws.Range("K10").Value = 5.2
$ws.Range("M10").Value = 1.03
$ws.Range("N10").Value = 3.45
$ws.Range("Q10").Value = 1.49
$ws.Range("R10").Value = 1.61
$ws.Range("S10").Value = 2.08
$ws.Range("T10").Value = 1.57
$ws.Range("U10").Value = 2.4
$ws.Range("V10").Value = 1.25
$ws.Range("W10").Value = 2.4
$ws.Range("X10").Value = 27
$ws.Range("Y10").Value = 26
$ws.Range("Z10").Value = 1000
$ws.Range("AB11").Value = 7.4
$ws.Range("F11").Value = 1.63
$ws.Range("G11").Value = 1.88
$ws.Range("H11").Value = 5.6
$ws.Range("I11").Value = 9
$ws.Range("J11").Value = 3.05
$ws.Range("K11").Value = 4.5
$ws.Range("L11").Value = 1.43
$ws.Range("N11").Value = 2.72
$ws.Range("O11").Value = 1.48
$ws.Range("P11").Value = 1.58
$ws.Range("Q11").Value = 2.36
$ws.Range("R11").Value = 1.21
$ws.Range("S11").Value = 4.2
$ws.Range("T11").Value = 2.06
$ws.Range("U11").Value = 1.64
$ws.Range("V11").Value = 1.12
$ws.Range("W11").Value = 2.12
$ws.Range("I12").Value = 3.1
$ws.Range("O12").Value = 1.5
$ws.Range("Q12").Value = 2.48
$ws.Range("V12").Value = 1.48
$ws.Range("Y12").Value = 10
$ws.Range("AB13").Value = 7.6
$ws.Range("AC13").Value = 8
$ws.Range("AD13").Value = 980
$ws.Range("AF13").Value = 12.5
$ws.Range("AG13").Value = 12
$ws.Range("AH13").Value = 28
$ws.Range("AL13").Value = 55
$ws.Range("AN13").Value = 980
$ws.Range("F13").Value = 1.92
$ws.Range("G13").Value = 2.06
$ws.Range("J13").Value = 3.4
$ws.Range("K13").Value = 3.6
$ws.Range("P13").Value = 1.6
$ws.Range("R13").Value = 1.23
$ws.Range("V13").Value = 1.25
$ws.Range("W13").Value = 1.95
$ws.Range("Y13").Value = 980
$ws.Range("AC14").Value = 7.8
$ws.Range("AI14").Value = 95
$ws.Range("AL14").Value = 85
$ws.Range("AO14").Value = 100
$ws.Range("F14").Value = 2.24
$ws.Range("G14").Value = 2.48
$ws.Range("H14").Value = 3.55
$ws.Range("I14").Value = 3.95
$ws.Range("J14").Value = 3.05
$ws.Range("K14").Value = 3.35
$ws.Range("M14").Value = 1.1
$ws.Range("P14").Value = 1.57
$ws.Range("S14").Value = 4.4
$ws.Range("T14").Value = 2.02
$ws.Range("U14").Value = 1.78
$ws.Range("V14").Value = 1.33
$ws.Range("W14").Value = 1.67
$ws.Range("AA15").Value = 1000
$ws.Range("AN15").Value = 1000
$ws.Range("AO15").Value = 90
$ws.Range("G15").Value = 2.46
$ws.Range("H15").Value = 3.6
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 3.1
$ws.Range("K15").Value = 3.2
$ws.Range("N15").Value = 2.6
$ws.Range("O15").Value = 1.54
$ws.Range("U15").Value = 1.74
$ws.Range("V15").Value = 1.33
$ws.Range("W15").Value = 1.68
$ws.Range("Z15").Value = 26
